$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing rows (4-6) that only held stray SONG_ID values.
$ws.Rows("4:6").Delete()

# Row 2: song "sw2op"
$ws.Range("A2").Value = "sw2op"
$ws.Range("C2").Value = "まいにちがドンダフル"
$ws.Range("D2").Value = "「太鼓の達人 ドンダフルフェスティバル」テーマソング"
$ws.Range("B2").Value = 5
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 2

# Row 3: song "yumeut" (fixes the balloon count that had been copied from the wrong row)
$ws.Range("C3").Value = "夢うつつカタルシス"
$ws.Range("A3").Value = "yumeut"
$ws.Range("D3").Value = "大木奏弥(BNSI) feat. 愛原圭織"
$ws.Range("B3").Value = 5
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 4

# Match the saved selection state from the authored workbook.
$ws.Range("J11").Select()
